# Rename the original (and only) worksheet to "data", then append three
# new blank worksheets (Sheet2, Sheet3, Sheet4) after it, and finally move
# the selection on the "data" sheet from B5 to D13 while keeping "data"
# the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "data"

$sheet2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$sheet2.Name = "Sheet2"

$sheet3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet2)
$sheet3.Name = "Sheet3"

$sheet4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet3)
$sheet4.Name = "Sheet4"

$ws.Activate() | Out-Null
$ws.Range("D13").Select() | Out-Null
